$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.807.61'
$ws.Range('E2').Value = '  -0.95%  '

$ws.Range('D3').Value = '1.942.29'
$ws.Range('E3').Value = '  -0.68%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.00'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.53%  '

$ws.Range('E6').Value = '  +0.10%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4889'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.05%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2949'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.59%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06884'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.55%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.48'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.29%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '106.29'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.05%  '

$ws.Range('D12').Value = '1.943.51'
$ws.Range('E12').Value = '  -0.51%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07726'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.16%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.362'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.86%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6981'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.02%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '276.57'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.13%  '

$ws.Range('D17').Value = '30.816.36'
$ws.Range('E17').Value = '  -0.71%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007725'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.59%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.12'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.77%  '

$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.06%  '

$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.193.23'
$ws.Range('E21').Value = '  -0.30%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.482'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.03%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.08%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.542'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.47%  '

$ws.Range('E25').Value = '  -2.00%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.31'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.29%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.68'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.07%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.166'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.02%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1046'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.48%  '

$ws.Range('E30').Value = '  -3.32%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.570'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.89%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.556'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.45%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.378'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.89%  '

$ws.Range('E34').Value = '  -2.77%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7542'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.41%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.160'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.52%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9998'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.12%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.733'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.12%  '

$ws.Range('E39').Value = '  -2.50%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.661'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.90%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.532'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.47%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '78.04'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.34%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.103'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.10%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9077'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.97%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '108.36'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.27%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4402'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.47%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9991'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.05%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.766'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.65%  '

$ws.Range('D49').Value = '1.002.14'
$ws.Range('E49').Value = '  +0.71%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1247'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.76%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.286'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.15%  '
